$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Refresh the "time_taken" column with the new panel-query run timestamps.
$ws.Range("F2").Value = "2021-10-05 14:19:38.697283"
$ws.Range("F3").Value = "2021-10-05 14:19:38.697292"
$ws.Range("F4").Value = "2021-10-05 14:19:38.697296"
$ws.Range("F5").Value = "2021-10-05 14:19:38.697298"
$ws.Range("F6").Value = "2021-10-05 14:19:38.697302"
$ws.Range("F7").Value = "2021-10-05 14:19:38.697304"
$ws.Range("F8").Value = "2021-10-05 14:19:38.697307"
$ws.Range("F9").Value = "2021-10-05 14:19:38.697310"
$ws.Range("F10").Value = "2021-10-05 14:19:38.697313"
$ws.Range("F11").Value = "2021-10-05 14:19:38.697316"
$ws.Range("F12").Value = "2021-10-05 14:19:38.697319"
$ws.Range("F13").Value = "2021-10-05 14:19:38.697321"
$ws.Range("F14").Value = "2021-10-05 14:19:38.697324"

# Add the new "metadata" tab right after "data".
$meta = $wb.Worksheets.Add($null, $ws)
$meta.Name = "metadata"

# Reuse the bold/bordered header style already used on the "data" sheet.
$ws.Range("B1:F1").Copy()
$meta.Range("B1:F1").PasteSpecial(-4122)
$ws.Range("B1").Copy()
$meta.Range("G1").PasteSpecial(-4122)

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Reuse the row-index style from the "data" sheet for A2.
$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)
$meta.Range("A2").Value = 0

$meta.Range("B2").Value = "Colorectal cancer pertinent cancer susceptibility"
$meta.Range("C2").Value = 244

# "1.0" must stay textual (matches the source export), not become the number 1.
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "1.0"

$meta.Range("E2").Value = "2017-11-05T02:37:20.290684Z"
$meta.Range("F2").Value = "2021-10-05 14:19:38.693572"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/244/?format=json"

# Keep "data" as the active/selected sheet, matching the original bookView.
$ws.Activate()
$ws.Range("A1").Select()
